$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 2 -----
$ws.Range("E2").Value = 23.40000000000022

$ws.Range("H2").Value = [double]"2.166288828536891e-16"
$ws.Range("I2").Value = 0.3453319618501958

$ws.Range("K2").Value = 49.00137610465953
$ws.Range("L2").Value = "[38.31771489841133, 59.68503731090774]"

$ws.Range("O2").Value = 1.867974010242579
$ws.Range("P2").Value = "[1.6415529180919632, 2.0943951023931957]"

$ws.Range("S2").Value = 68.56513316854011
$ws.Range("T2").Value = "[62.41390728030426, 74.71635905677597]"

$ws.Range("W2").Value = 16.4432432432434
$ws.Range("X2").Value = 15.60000000000015
$ws.Range("Y2").Value = 17.28648648648665

# ----- Row 3 -----
$ws.Range("B3").Value = 1

$ws.Range("E3").Value = 22.91000000000014

$ws.Range("G3").Value = [double]"1.110223024625157e-16"
$ws.Range("H3").Value = [double]"2.166288828536891e-16"

$ws.Range("K3").Value = 58.90988150023923
$ws.Range("L3").Value = "[43.96084979893658, 73.85891320154188]"

$ws.Range("M3").Value = [double]"4.771738559838923e-13"
$ws.Range("N3").Value = [double]"4.771738559838923e-13"

$ws.Range("O3").Value = 1.037763339023655
$ws.Range("P3").Value = "[0.761026448617347, 1.3145002294299637]"

$ws.Range("Q3").Value = [double]"4.411360166045597e-12"
$ws.Range("R3").Value = [double]"4.411360166045597e-12"

$ws.Range("S3").Value = 68.18068274826216
$ws.Range("T3").Value = "[60.286301000950594, 76.07506449557371]"

$ws.Range("W3").Value = 19.12606606606618
$ws.Range("X3").Value = 18.11701701701713
$ws.Range("Y3").Value = 20.13511511511524
